$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Albert Einstein): replace Biography text (E3) with the expanded paragraph
$ws.Range("E3").Value = 'Einstein was not able to speak until he was 4 years old and his teachers said "he would never amount to much". Einstein’s parents wanted him to pursue a career in electrical engineering. Without a high school diploma, Einstein applied to the Polytechnic Institute at Zurich, Switzerland and failed the entrance examination although he got exceptional marks in the mathematics and physics sections. Einstein described himself as “a conscientious but unassuming young man who had acquired his meagre store of pertinent knowledge of the essentials through self study." Einstein basically taught himself math and physics through books and self-study and skipped lectures and stopped going to classes that didn''t interest him. When he was working on the Theory of Relativity, he was working as an Assistant Examiner at a patent office. And then of course, he became one of the greatest scientists of all time. '

# Row 22
$ws.Range("A22").Value = 'Steve Jobs'
$ws.Range("C22").Value = 'Depression, Mood Disorders, Lack of Formal Education'
$ws.Range("E22").Value = 'Steve Jobs was a famous college-dropout who was famous for his severe mood swings. At 30, he was left devastated and depressed after being fired from Apple, the very company he started. When Jobs looked back at this setback, he said this setback pushed him to "One of the most creative periods of his life" where he eventually founded Pixar and NextStep. His legacy is what allowed me to create this iOS app for you ;)'

# Row 23
$ws.Range("A23").Value = 'Walt Disney'
$ws.Range("E23").Value = 'Disney was once fired from a newspaper company for "lacking imagination" and "having no original thoughts". 
Goes to show that other people are a terrible judge of your internal potential ;) Probably because they don''t know that ''what is essential is invisible to the eye'' '

# Row 24
$ws.Range("A24").Value = 'Honda'

# Row 25
$ws.Range("A25").Value = 'Michael Jordan'
$ws.Range("C25").Value = 'Late Bloomer'
$ws.Range("E25").Value = 'He was cut from his HS basketball team and apparently went home, locked himself in his room and cried. 
And those tears then propelled him into becoming the greatest basketball player of all-time.'

# Row 26
$ws.Range("A26").Value = 'Paulo Coelho'
$ws.Range("C26").Value = 'Rejection'
$ws.Range("E26").Value = '"The Alchemist" was rejected over 200 times before it went on to sell 75 Million cpies. '

# Row 27
$ws.Range("A27").Value = 'Oprah Winfrey'
$ws.Range("C27").Value = 'Childhood abuse'
$ws.Range("E27").Value = 'Oprah was sexualy abused by male relatives of her family during adolescence. 
Today, she is one of the most beloved TV show talk hosts, actress, publisher, producer and philanthropist.'

# Row 28
$ws.Range("A28").Value = 'Richard Pryor'
$ws.Range("C28").Value = 'Childhood abuse, Stage Fright'

# Row 29
$ws.Range("A29").Value = 'Rodney Dangerfield'
$ws.Range("C29").Value = 'Late Bloomer'

# Row 30
$ws.Range("A30").Value = 'Jim Carrey'
$ws.Range("C30").Value = 'Poverty, Depression'

# Row height adjustments to match wrapped content
$ws.Rows.Item(22).RowHeight = 248
$ws.Rows.Item(23).RowHeight = 155
$ws.Rows.Item(25).RowHeight = 124
$ws.Rows.Item(26).RowHeight = 62
$ws.Rows.Item(27).RowHeight = 124
$ws.Rows.Item(28).RowHeight = 62

# Update view state: frozen-pane scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 6
$ws.Range("G28").Select()
